$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Nikolas, Timeo, Malkovan"

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "Alecxperdu, Nelson, DTP, DramaPanda, Maitredoudou, FooD_Flo"

$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "Alecxperdu_S2, Gobou, Polo, Skowa, Mark_S2"

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "Krak, Nwog, Ethan, Jehovah, Alexadventure, Step, Feiik, Tiff, Natoxe, Bilal, Flau_S4, Didine, Bilel_S5, TimeoGnc"

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "Lili, Samson, SorciShoot, Louan, Xo, Napoleon, Cyrf_S4, Espoir Perdu, And_S5, Major Chris, Nikolas_S6, Nemocca, Aurel, Theo, Chatoon, Clement, Mizuki, Armand, Manu"

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "Chreet, Joshua, Gririsu, Hugo, Thib, Julien, Mickey_S8, Sunka"

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Chreet_S5"

$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Kamiga, Z4tix, Mark, Gwendal, Clem fair play, Amaury, Schweppes"

$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "Nicovid, Ladoly, Deku, Sayo, Line_S5, Clem fair play_S5"

$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Blgham"

$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "The soulless, FooD_Flo_S5, Julien_S7"

$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "Black Lolo, Eneko, Enze, Thyx, Doggydog_S6, Tim, Luc"

$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "Schweppes_S6, Ethan_S6, Xori, Gigi, Spider, Enze_S8"

$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "Kamiga_S5, GiulfeuYT_S6, Alan_S7, Sy_boulette, Reiko, Ju, Mielle"

$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = "Yuu, Nounours, Saminette, Nicovid_S6, Cyrf_S6, Blacks Star, Kamiga_S9"

$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = "Saucisson, Ranzyo_S5, Killian, Waikato, Mickey, Jilink, Pilouche, Alan_S9"

$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = "Nelson_S4, Alex., Skowa_S4, Raphi, Xiao_S5, Manta, Dahmi1 Arti, Gobou_S6, Gwendal_S7, Enze_S7, Lilian, Mehdiiii, Corentin"

$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = "Bylost, Line, HiYoucef, Saucisson_S5, Malkovan_S5, Quentin, Alex"

$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = "Crypto, Yatho, Vah Balress, Mistimat, Guigui_S9, Mtking"

$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = "Flau"

$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = "Bilel_S4, Bylost_S4, Thib_S5, Grenzo, Chopa, Cha, Z4tix_S9"

$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = "And, Xiao, Gririsu_S5, Nath_S6, Clovis, Thynael, Cosmos, Piiskoo"

$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = "Cyrf, Sneus, Malkovan_S3, Julien_S5, Timeo_S5, Black Lolo_S5, Nikolas_S5, Amaury_S5, Kwinn"

$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = "Nath, Ethan_S3, Ranzyo, Guigui, Angel, Dragon"

$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = "Bilel, Joshua_S3, GiulfeuYT, Doggydog, Alan, Mark_S5, Maitredoudou_S5"
